# "With held reconcilation fee" - refresh the latest automated test run's
# reconciliation data on the 1099MISCdata sheet, and tidy up the active
# window/selection state left over from that run.

$wb = $excel.ActiveWorkbook

# --- 1099MISCdata: latest reconciliation test-run values ---------------
$ws = $wb.Worksheets.Item("1099MISCdata")
$ws.Activate()

$ws.Range("B3").Value = "Test0310202165905"   # Recipient BusinessName (test run id)
$ws.Range("B4").Value = "Test0310202165232"   # Payer last name (test run id)

# Reference number is digits-only, so force text so it doesn't get coerced
# into a number cell.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "201204749"           # Reference number

$ws.Range("B8").Value = "03/15/2021"          # SchedlueDate
$ws.Range("B9").Value = "96-8531926"          # Recipient EIN
$ws.Range("B16").Value = 162021761            # Payer EIN number

# Leave the selection on the freshly-updated payer-last-name cell, with the
# view scrolled back to the top of the sheet.
$ws.Range("A1").Select()
$ws.Range("B4").Select()

# --- Addpayerrandombusinessein: drop the stale selection ---------------
$ws2 = $wb.Worksheets.Item("Addpayerrandombusinessein")
$ws2.Activate()
$ws2.Range("A1").Select()

# restore original active sheet
$ws.Activate()

# --- Application window size -------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 16920
$win.Height = 8220
